# Generate Report for Handoff
# Updates the "Latest Handoff Datetime"/"Latest Handoff Date" timestamps for the
# files that were (re-)handed off: 5bf95739..., 0136af66..., 7a044612...,
# 7a9f6efc..., 9ece3265..., ade475ab..., bc3b1be2..., dbdf809c...

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest Handoff Date" column (D) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D10").Value = "2016-26-18 20:26:57"
$wsOverview.Range("D13").Value = "2016-26-18 20:26:57"

# --- zh-cn sheet: "Latest Handoff Datetime" column (E) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E7").Value = "2016-03-18 20:26:54"
$wsZhCn.Range("E10").Value = "2016-03-18 20:26:54"
$wsZhCn.Range("E11").Value = "2016-03-18 20:26:54"
$wsZhCn.Range("E12").Value = "2016-03-18 20:26:54"
$wsZhCn.Range("E13").Value = "2016-03-18 20:26:54"
$wsZhCn.Range("E14").Value = "2016-03-18 20:26:54"
$wsZhCn.Range("E15").Value = "2016-03-18 20:26:54"
$wsZhCn.Range("E16").Value = "2016-03-18 20:26:54"

# --- de-de sheet: "Latest Handoff Datetime" column (E) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E7").Value = "2016-03-18 20:26:57"
$wsDeDe.Range("E10").Value = "2016-03-18 20:26:57"
$wsDeDe.Range("E11").Value = "2016-03-18 20:26:57"
$wsDeDe.Range("E12").Value = "2016-03-18 20:26:57"
$wsDeDe.Range("E13").Value = "2016-03-18 20:26:57"
$wsDeDe.Range("E14").Value = "2016-03-18 20:26:57"
$wsDeDe.Range("E15").Value = "2016-03-18 20:26:57"
$wsDeDe.Range("E16").Value = "2016-03-18 20:26:57"
